$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new header cell F1 ("time_taken"), matching the style of the
# existing header cells (e.g. E1: bold, bordered, centered/top-aligned).
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Per-row timestamps captured when each panel entry was processed.
$timeTaken = @{
    2  = "2021-10-05 13:40:16.445329"
    3  = "2021-10-05 13:40:16.445342"
    4  = "2021-10-05 13:40:16.445346"
    5  = "2021-10-05 13:40:16.445349"
    6  = "2021-10-05 13:40:16.445353"
    7  = "2021-10-05 13:40:16.445356"
    8  = "2021-10-05 13:40:16.445359"
    9  = "2021-10-05 13:40:16.445362"
    10 = "2021-10-05 13:40:16.445366"
    11 = "2021-10-05 13:40:16.445369"
    12 = "2021-10-05 13:40:16.445372"
    13 = "2021-10-05 13:40:16.445375"
    14 = "2021-10-05 13:40:16.445378"
    15 = "2021-10-05 13:40:16.445381"
    16 = "2021-10-05 13:40:16.445384"
    17 = "2021-10-05 13:40:16.445387"
    18 = "2021-10-05 13:40:16.445390"
    19 = "2021-10-05 13:40:16.445393"
    20 = "2021-10-05 13:40:16.445397"
    21 = "2021-10-05 13:40:16.445399"
    22 = "2021-10-05 13:40:16.445402"
    23 = "2021-10-05 13:40:16.445405"
    24 = "2021-10-05 13:40:16.445408"
    25 = "2021-10-05 13:40:16.445411"
    26 = "2021-10-05 13:40:16.445415"
}

foreach ($row in 2..26) {
    $ws.Cells.Item($row, 6).Value = $timeTaken[$row]
}
